# Append 4 new "rcp45" run-summary rows to the "Run Info" sheet, mirroring
# the formatting of the existing rows by copying the last data row down and
# then overwriting the values that differ.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run Info")

$srcRow = 205
$startRow = 206

# Columns: A Date, B RCP, C E, D everyx, E workers, F run time,
# G 85-2010 bleaching, H sBleach 1, I cBleach 1, J sRecSeedMult 1,
# K cRecSeedMult 1, L cSeedThreshMult 1, M pMin, N pMax, O exponent, P divisor
$newRows = @(
    @(42962.527442129627, 1, 1,   6, 19.147612058121716, 5.0069930069930075),
    @(42962.527881944443, 1, 1,   6, 21.710316191785406, 5.0069930069930075),
    @(42962.530034722222, 1, 100, 6, 7.8395451405188386, 6.4685314685314683),
    @(42962.544351851851, 1, 100, 6, 26.86526018232394,  6.4685314685314683)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $destRow = $startRow + $i
    $data = $newRows[$i]

    # Copy formats + values from the template row, then fix up the cells
    # that actually differ for this run.
    $ws.Range("A${srcRow}:P${srcRow}").Copy($ws.Range("A${destRow}:P${destRow}"))

    $ws.Cells.Item($destRow, 1).Value = $data[0]   # Date
    $ws.Cells.Item($destRow, 3).Value = $data[1]   # E
    $ws.Cells.Item($destRow, 4).Value = $data[2]   # everyx
    $ws.Cells.Item($destRow, 5).Value = $data[3]   # workers
    $ws.Cells.Item($destRow, 6).Value = $data[4]   # run time
    $ws.Cells.Item($destRow, 7).Value = $data[5]   # 85-2010 bleaching
}

$lastRow = $startRow + $newRows.Count - 1
$ws.Range("A${lastRow}:P${lastRow}").Select() | Out-Null
